# Add a 4-wise comparison experiment:
#  - Solver sheet: append new simulation-result rows (19-22, 24-27; rows 18
#    and 23 are left blank as separators, matching the existing workbook
#    convention of blank rows between parameter blocks) and grow Table1 to
#    cover them.
#  - Data sheet: fill in the previously-empty "Numb agents"/"Time(sec)"
#    columns (E/F) for rows 11-16.
#  - Move the active selection from the Data sheet to the Solver sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Solver")
$ws4 = $wb.Worksheets.Item("Data")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

# New Solver rows keyed by row number -> full A:P row values.
# Rows 18 and 23 are intentionally skipped (kept blank), mirroring the
# blank separator rows already used elsewhere in this workbook (e.g. the
# "New Cost" sheet's Table3).
$rowsData = @{
    19 = @(10, "Random",   "Integer Prog", "No",  1220.3493382930701, 2.0501657930400001, 1.80033615631276,   1.7194241585314101, 2.0325430464415102, 2.2803153108354199, 1.696228990589,     2.00750640695327,   1.76022477113234,   1.46777654373842,   3.0045180521204302, 2.7327844937752701)
    20 = @(10, "Random",   "Integer Prog", "Yes", 1220.3493382930701, 2.0501233599200002, 1.80033615631276,   1.7194241585314101, 2.0325430464415102, 2.2803153108354199, 1.69621524332908,   2.0075041315309901, 1.75985905837951,   1.4677339479690901, 3.0045180521204302, 2.7327844937752701)
    21 = @(10, "Random",   "Lottery",      "No",  1220.3493382930701, 1.97745122474,      1.7511793390152,    1.6267033808568501, 1.95912559500413,   2.2346833033150699, 1.6816245289632099, 1.94344637256716,   1.7180614688439999, 1.4441184865127199, 2.86398405855922,   2.5515857137933602)
    22 = @(10, "Random",   "Lottery",      "Yes", 1220.3493382930701, 1.97600739867,      1.7511793390152,    1.6267033808568501, 1.95912559500413,   2.2346833033150699, 1.6816245289632099, 1.93562900734996,   1.7180614688439999, 1.4441184865127199, 2.86398405855922,   2.5449648183154201)
    24 = @(10, "Location", "Integer Prog", "No",  1236.0917141437501, 2.4165244544100002, 1.17274621812683,   1.47817270300976,   1.20159567625994,   3.2420457974741699, 2.0143376002501698, 1.9534244929914999, 9.2755471800342306, 1.47552103335132,   1.20640070875853,   1.14545313381587)
    25 = @(10, "Location", "Integer Prog", "Yes", 1236.0917141437501, 2.4165241577300001, 1.17274616599535,   1.47817270300976,   1.20159567625994,   3.2420457974741699, 2.01433766502227,   1.9534214872208899, 9.2755471562492602, 1.47552108345765,   1.20640070875853,   1.14545313381587)
    26 = @(10, "Location", "Lottery",      "No",  1236.0917141437501, 2.2672542904599999, 1.1726595112798299, 1.4606033441969899, 1.2014641699966999, 3.0994407377636399, 1.96882958188343,   1.8822722061053001, 8.0675666415269305, 1.4724767692734899, 1.2017770237582699, 1.1454529188434099)
    27 = @(10, "Location", "Lottery",      "Yes", 1236.0917141437501, 2.2659206216699999, 1.1726595112798299, 1.4466640895373899, 1.1927889491158401, 3.1000708665018499, 1.95225661151936,   1.9088426449232001, 8.0675666415269305, 1.4711277991582099, 1.2017770237582699, 1.14545207941674)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws1.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# Grow Table1 so its range/autofilter covers the newly-added rows (blank
# rows 18 & 23 stay inside the table range, same as the blank rows already
# embedded inside Table3 on the "New Cost" sheet).
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:P27"))

# Fill in the previously blank "Numb agents" (E) / "Time(sec)" (F) columns
# on the Data sheet for rows 11-16.
$data4 = @{
    11 = @(801,  4)
    12 = @(1014, 7)
    13 = @(1179, 6)
    14 = @(424,  2)
    15 = @(319,  3)
    16 = @(319,  3)
}

foreach ($r in $data4.Keys) {
    $vals = $data4[$r]
    $ws4.Range("E$r").Value = $vals[0]
    $ws4.Range("F$r").Value = $vals[1]
}

# Update selections: Data sheet keeps a (non-active) selection at J9, while
# the Solver sheet becomes the active tab with F22 selected.
$ws4.Activate()
$ws4.Range("J9").Select()

$ws1.Activate()
$ws1.Range("F22").Select()
